# Actualización automática 2025-06-23 17:05:09
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("M46").Value = 171.19
$wsVentasPorGrupo.Range("D51").Value = 380.16
$wsVentasPorGrupo.Range("D54").Value = "7 de 52"
$wsVentasPorGrupo.Range("M54").Value = "8 de 52"

# --- Sheet "VENTA MENSUAL" ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F46").Value = 171.19
$wsVentaMensual.Range("F51").Value = 380.16
$wsVentaMensual.Range("F54").Value = 47135.57

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimientoMensual = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3: 240X80 PORCELANATO
$wsCumplimientoMensual.Range("D3").Value = 8930.139999999999
$wsCumplimientoMensual.Range("E3").Value = 18526.8676
$wsCumplimientoMensual.Range("F3").Value = 0.3252408321437038

# Row 16: PORCELANATO
$wsCumplimientoMensual.Range("D16").Value = 9836.41
$wsCumplimientoMensual.Range("E16").Value = 22905.04
$wsCumplimientoMensual.Range("F16").Value = 0.3004268289889421

# Row 19: TOTAL
$wsCumplimientoMensual.Range("D19").Value = 47135.57
$wsCumplimientoMensual.Range("E19").Value = 47311.87064517915
$wsCumplimientoMensual.Range("F19").Value = 0.4990666732524734
